$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 0.09866275701177361

# Row 3
$ws.Range("B3").Value = 0.07117535186157629
$ws.Range("H3").Value = 0.1698381088733499

# Row 4
$ws.Range("B4").Value = 0.0875375577432094
$ws.Range("H4").Value = 0.186200314754983

# Row 5
$ws.Range("B5").Value = 0.08929271586081422
$ws.Range("H5").Value = 0.1879554728725878

# Row 6
$ws.Range("B6").Value = 0.1297238721402159
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("H6").Value = 0.2283866291519895

# Row 7
$ws.Range("B7").Value = 0.01991737152930807
$ws.Range("C7").Value = 0.004724198655420279
$ws.Range("D7").Value = 0.8523608007022909
$ws.Range("E7").Value = 0.02596587967591011
$ws.Range("F7").Value = 0.01063935756979954
$ws.Range("G7").Value = 0.02919538548881713
$ws.Range("H7").Value = 0.1185801285410817

# Row 8
$ws.Range("B8").Value = 0.01202422563404924
$ws.Range("C8").Value = 0.003079965097511697
$ws.Range("D8").Value = 0.581425895044372
$ws.Range("E8").Value = 0.0247495421370266
$ws.Range("F8").Value = 0.005977679743849167
$ws.Range("G8").Value = 0.01807077152424898
$ws.Range("H8").Value = 0.1106869826458229

# Row 9
$ws.Range("B9").Value = 0.02075623521033395
$ws.Range("C9").Value = 0.005028289645217988
$ws.Range("D9").Value = 0.6920075986725491
$ws.Range("E9").Value = 0.0119419471126582
$ws.Range("F9").Value = 0.0108343923415104
$ws.Range("G9").Value = 0.0306780780791569
$ws.Range("H9").Value = 0.1194189922221076

# Row 10
$ws.Range("B10").Value = 0.01876962931216891
$ws.Range("C10").Value = 0.003990427111981799
$ws.Range("D10").Value = 0.563613201886343
$ws.Range("E10").Value = 0.01197788420240421
$ws.Range("F10").Value = 0.01092969134454767
$ws.Range("G10").Value = 0.02660956727978983
$ws.Range("H10").Value = 0.1174323863239425

# Row 11
$ws.Range("B11").Value = 0.02665439086576339
$ws.Range("H11").Value = 0.125317147877537

# Row 12
$ws.Range("B12").Value = 0.04444008523698022
$ws.Range("H12").Value = 0.1431028422487538

# Row 13
$ws.Range("B13").Value = 0.05284586330793421
$ws.Range("H13").Value = 0.1515086203197078

# Row 14
$ws.Range("B14").Value = 0.05929096156625544
$ws.Range("H14").Value = 0.1579537185780291

# Row 15
$ws.Range("B15").Value = 0.0647307532444439
$ws.Range("H15").Value = 0.1633935102562175

# Row 16
$ws.Range("B16").Value = 0.06698157313446509
$ws.Range("H16").Value = 0.1656443301462387

# Row 17
$ws.Range("B17").Value = 0.07122168424452854
$ws.Range("H17").Value = 0.1698844412563021

# Row 18
$ws.Range("B18").Value = -0.09866275701177361
$ws.Range("C18").Value = 0.01293457410127628
$ws.Range("D18").Value = -15.98674484932656
$ws.Range("E18").Value = 0.04324973861679528
$ws.Range("F18").Value = -0.124071371995926
$ws.Range("G18").Value = -0.0732541420276209

# Row 19
$ws.Range("B19").Value = 0.07367202488660815
$ws.Range("H19").Value = 0.1723347818983818

# Row 20
$ws.Range("B20").Value = 0.07328462749031499
$ws.Range("H20").Value = 0.1719473845020886

# Row 21
$ws.Range("B21").Value = 0.0774210796666675
$ws.Range("H21").Value = 0.1760838366784411

# Row 22
$ws.Range("B22").Value = 0.08081373114732876
$ws.Range("C22").Value = 0.009102482634710026
$ws.Range("D22").Value = 16.21609148284125
$ws.Range("E22").Value = 0.03801489172310542
$ws.Range("F22").Value = 0.06293598043017892
$ws.Range("G22").Value = 0.09869148186447806
$ws.Range("H22").Value = 0.1794764881591024

# Row 23
$ws.Range("B23").Value = 0.08075664033354045
$ws.Range("C23").Value = 0.009093407507047017
$ws.Range("D23").Value = 15.67383018286934
$ws.Range("E23").Value = 0.04648444410622087
$ws.Range("F23").Value = 0.0628853429116049
$ws.Range("G23").Value = 0.09862793775547635
$ws.Range("H23").Value = 0.1794193973453141

# Row 24
$ws.Range("B24").Value = 0.0765923500713548
$ws.Range("C24").Value = 0.009400130650966678
$ws.Range("D24").Value = 13.9738056122143
$ws.Range("E24").Value = 0.04932007785943611
$ws.Range("F24").Value = 0.05811731122006756
$ws.Range("G24").Value = 0.09506738892264226
$ws.Range("H24").Value = 0.1752551070831284

# Row 25
$ws.Range("B25").Value = 0.0775071489822161
$ws.Range("C25").Value = 0.010352278973985
$ws.Range("D25").Value = 12.536241077218
$ws.Range("E25").Value = 0.06068493445085101
$ws.Range("F25").Value = 0.05713154062896972
$ws.Range("G25").Value = 0.09788275733546273
$ws.Range("H25").Value = 0.1761699059939897

# Row 26
$ws.Range("B26").Value = 0.07571473007823297
$ws.Range("C26").Value = 0.009773697408846009
$ws.Range("D26").Value = 11.32295039087089
$ws.Range("E26").Value = 0.07023439580074119
$ws.Range("F26").Value = 0.05650815938003449
$ws.Range("G26").Value = 0.09492130077643143
$ws.Range("H26").Value = 0.1743774870900066

# Row 27
$ws.Range("B27").Value = 0.07637824796035497
$ws.Range("C27").Value = 0.009739343173230108
$ws.Range("D27").Value = 10.95541747595428
$ws.Range("E27").Value = 0.07410456790748245
$ws.Range("F27").Value = 0.05723155328387933
$ws.Range("G27").Value = 0.09552494263683049
$ws.Range("H27").Value = 0.1750410049721286

# Row 28
$ws.Range("B28").Value = 0.07424768192339194
$ws.Range("C28").Value = 0.009322504667293849
$ws.Range("D28").Value = 10.50761822850089
$ws.Range("E28").Value = 0.10448923277448
$ws.Range("F28").Value = 0.05592103738459381
$ws.Range("G28").Value = 0.09257432646218938
$ws.Range("H28").Value = 0.1729104389351656

# Row 29
$ws.Range("B29").Value = 0.0227032714215708
$ws.Range("C29").Value = 0.003812182640184231
$ws.Range("D29").Value = 1.010503939714208
$ws.Range("E29").Value = 0.001264054602835931
$ws.Range("F29").Value = 0.01520386516299619
$ws.Range("G29").Value = 0.03020267768014518
$ws.Range("H29").Value = 0.1213660284333444
